$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1658.091
$ws.Range("I9").Value = 138
$ws.Range("K9").Value = 138
$ws.Range("M9").Value = 31
$ws.Range("H16").Value = 11163.167
$ws.Range("I16").Value = 7991
$ws.Range("J16").Value = 11797.6
$ws.Range("K16").Value = 7991
$ws.Range("L16").Value = 11797.6
$ws.Range("M16").Value = -7761
$ws.Range("N16").Value = -12257.6
$ws.Range("H41").Value = 4249.5
$ws.Range("I41").Value = 4499
$ws.Range("K41").Value = 4499
$ws.Range("M41").Value = -4059
$ws.Range("H51").Value = 12987.3
$ws.Range("J51").Value = 3600.25
$ws.Range("L51").Value = 3600.25
$ws.Range("N51").Value = -4568.25
$ws.Range("H69").Value = 22556.277
$ws.Range("I69").Value = 22556.277
$ws.Range("K69").Value = 67668.83099999999
$ws.Range("M69").Value = -66794.83099999999
$ws.Range("H72").Value = 22556.277
$ws.Range("I72").Value = 22556.277
$ws.Range("K72").Value = 203006.493
$ws.Range("M72").Value = -198638.493
$ws.Range("H80").Value = 989.3333
$ws.Range("I80").Value = 274.45456
$ws.Range("J80").Value = 1594.2307
$ws.Range("K80").Value = 823.36368
$ws.Range("L80").Value = 4782.6921
$ws.Range("M80").Value = 174.63632
$ws.Range("N80").Value = -6778.6921
$ws.Range("H83").Value = 989.3333
$ws.Range("I83").Value = 274.45456
$ws.Range("J83").Value = 1594.2307
$ws.Range("K83").Value = 2470.09104
$ws.Range("L83").Value = 14348.0763
$ws.Range("M83").Value = 2521.90896
$ws.Range("N83").Value = -24332.0763
$ws.Range("H86").Value = 3490.348
$ws.Range("I86").Value = 2842.7144
$ws.Range("K86").Value = 2842.7144
$ws.Range("M86").Value = -1719.7144
$ws.Range("H89").Value = 3490.348
$ws.Range("I89").Value = 2842.7144
$ws.Range("K89").Value = 14213.572
$ws.Range("M89").Value = -8597.572
$ws.Range("H96").Value = 3642
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -17746
$ws.Range("H98").Value = 4034.842
$ws.Range("I98").Value = 4197.6875
$ws.Range("K98").Value = 4197.6875
$ws.Range("M98").Value = -2699.6875
$ws.Range("H103").Value = 589.5405
$ws.Range("I103").Value = 582.3214
$ws.Range("J103").Value = 612
$ws.Range("K103").Value = 1746.9642
$ws.Range("L103").Value = 1836
$ws.Range("M103").Value = -1160.9642
$ws.Range("N103").Value = -3008
$ws.Range("H122").Value = 4034.842
$ws.Range("I122").Value = 4197.6875
$ws.Range("K122").Value = 12593.0625
$ws.Range("M122").Value = -10143.0625
$ws.Range("H132").Value = 1570.186
$ws.Range("I132").Value = 1379.1842
$ws.Range("J132").Value = 3021.8
$ws.Range("K132").Value = 4137.5526
$ws.Range("L132").Value = 9065.400000000001
$ws.Range("M132").Value = -1607.5526
$ws.Range("N132").Value = -14125.4
$ws.Range("H137").Value = 1552.409
$ws.Range("I137").Value = 1040.4482
$ws.Range("K137").Value = 3121.3446
$ws.Range("M137").Value = -571.3446000000004
$ws.Range("H141").Value = 123453.875
$ws.Range("I141").Value = 140233
$ws.Range("K141").Value = 420699
$ws.Range("M141").Value = -415519

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1870.0741
$ws.Range("I2").Value = 1947.3889
$ws.Range("K2").Value = 1947.3889
$ws.Range("M2").Value = -1834.3889
$ws.Range("H32").Value = 3232008.8
$ws.Range("I32").Value = 3339442.8
$ws.Range("K32").Value = 3339442.8
$ws.Range("M32").Value = -3339155.8
$ws.Range("H45").Value = 2284.2307
$ws.Range("I45").Value = 1721.8948
$ws.Range("K45").Value = 1721.8948
$ws.Range("M45").Value = -1344.8948
$ws.Range("H61").Value = 3012
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H74").Value = 1907
$ws.Range("I74").Value = 1192.1111
$ws.Range("K74").Value = 1192.1111
$ws.Range("M74").Value = -318.1111000000001
$ws.Range("H77").Value = 1907
$ws.Range("I77").Value = 1192.1111
$ws.Range("K77").Value = 5960.5555
$ws.Range("M77").Value = -1592.5555
$ws.Range("H88").Value = 3845.6
$ws.Range("I88").Value = 1836
$ws.Range("K88").Value = 1836
$ws.Range("M88").Value = -1430
$ws.Range("H91").Value = 3845.6
$ws.Range("I91").Value = 1836
$ws.Range("K91").Value = 1836
$ws.Range("M91").Value = -432
$ws.Range("H97").Value = 9968.5
$ws.Range("I97").Value = 9968.5
$ws.Range("K97").Value = 9968.5
$ws.Range("M97").Value = -9472.5
$ws.Range("H110").Value = 3371.7144
$ws.Range("I110").Value = 1691.1875
$ws.Range("J110").Value = 8749.4
$ws.Range("K110").Value = 1691.1875
$ws.Range("L110").Value = 8749.4
$ws.Range("M110").Value = 353.8125
$ws.Range("N110").Value = -12839.4
$ws.Range("H116").Value = 1870.0741
$ws.Range("I116").Value = 1947.3889
$ws.Range("K116").Value = 1947.3889
$ws.Range("M116").Value = 346.6111000000001
$ws.Range("H122").Value = 3165
$ws.Range("I122").Value = 3182.3333
$ws.Range("J122").Value = 3113
$ws.Range("K122").Value = 9546.999899999999
$ws.Range("L122").Value = 9339
$ws.Range("M122").Value = -7096.999899999999
$ws.Range("N122").Value = -14239
$ws.Range("H132").Value = 2694.6191
$ws.Range("I132").Value = 1922
$ws.Range("J132").Value = 5978.25
$ws.Range("K132").Value = 5766
$ws.Range("L132").Value = 17934.75
$ws.Range("M132").Value = -3236
$ws.Range("N132").Value = -22994.75
$ws.Range("H136").Value = 3012
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1870.0741
$ws.Range("I3").Value = 1947.3889
$ws.Range("K3").Value = 1947.3889
$ws.Range("M3").Value = -1833.3889
$ws.Range("H31").Value = 13500
$ws.Range("I31").Value = 13500
$ws.Range("K31").Value = 13500
$ws.Range("M31").Value = -13248
$ws.Range("H55").Value = 49509.668
$ws.Range("J55").Value = 49509.668
$ws.Range("L55").Value = 49509.668
$ws.Range("N55").Value = -50055.668
$ws.Range("H86").Value = 7870.04
$ws.Range("I86").Value = 3519.6667
$ws.Range("K86").Value = 3519.6667
$ws.Range("M86").Value = -2396.6667
$ws.Range("H89").Value = 7870.04
$ws.Range("I89").Value = 3519.6667
$ws.Range("K89").Value = 17598.3335
$ws.Range("M89").Value = -11982.3335
$ws.Range("H105").Value = 2900
$ws.Range("I105").Value = 2900
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2900
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -1153
$ws.Range("H107").Value = 3380.5
$ws.Range("I107").Value = 2486.3333
$ws.Range("J107").Value = 5168.8335
$ws.Range("K107").Value = 2486.3333
$ws.Range("L107").Value = 5168.8335
$ws.Range("M107").Value = -566.3332999999998
$ws.Range("N107").Value = -9008.833500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 459.78946
$ws.Range("I7").Value = 112.870964
$ws.Range("J7").Value = 1996.1428
$ws.Range("K7").Value = 112.870964
$ws.Range("L7").Value = 1996.1428
$ws.Range("M7").Value = 0.1290359999999993
$ws.Range("N7").Value = -2222.1428
$ws.Range("H31").Value = 1722.3529
$ws.Range("I31").Value = 1346.9231
$ws.Range("J31").Value = 2942.5
$ws.Range("K31").Value = 1346.9231
$ws.Range("L31").Value = 2942.5
$ws.Range("M31").Value = -1051.9231
$ws.Range("N31").Value = -3532.5
$ws.Range("H34").Value = 1722.3529
$ws.Range("I34").Value = 1346.9231
$ws.Range("J34").Value = 2942.5
$ws.Range("K34").Value = 1346.9231
$ws.Range("L34").Value = 2942.5
$ws.Range("M34").Value = -1144.9231
$ws.Range("N34").Value = -3346.5
$ws.Range("H58").Value = 1651.2069
$ws.Range("I58").Value = 1287.619
$ws.Range("J58").Value = 2605.625
$ws.Range("K58").Value = 1287.619
$ws.Range("L58").Value = 2605.625
$ws.Range("M58").Value = -1084.619
$ws.Range("N58").Value = -3011.625
$ws.Range("H86").Value = 3933.3333
$ws.Range("J86").Value = 3999.5
$ws.Range("L86").Value = 3999.5
$ws.Range("N86").Value = -6245.5
$ws.Range("H89").Value = 3933.3333
$ws.Range("J89").Value = 3999.5
$ws.Range("L89").Value = 19997.5
$ws.Range("N89").Value = -31229.5
$ws.Range("H99").Value = 2226.818
$ws.Range("I99").Value = 1719.2858
$ws.Range("K99").Value = 1719.2858
$ws.Range("M99").Value = -221.2858000000001
$ws.Range("H126").Value = 2226.818
$ws.Range("I126").Value = 1719.2858
$ws.Range("K126").Value = 5157.857400000001
$ws.Range("M126").Value = -2687.857400000001
$ws.Range("H132").Value = 1758.6666
$ws.Range("I132").Value = 1678.2222
$ws.Range("K132").Value = 5034.6666
$ws.Range("M132").Value = -2504.6666
$ws.Range("H136").Value = 1651.2069
$ws.Range("I136").Value = 1287.619
$ws.Range("J136").Value = 2605.625
$ws.Range("K136").Value = 3862.857
$ws.Range("L136").Value = 7816.875
$ws.Range("M136").Value = -1312.857
$ws.Range("N136").Value = -12916.875
$ws.Range("H141").Value = 166954.77
$ws.Range("J141").Value = 179998.73
$ws.Range("L141").Value = 179998.73
$ws.Range("N141").Value = -190358.73

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 167278.33
$ws.Range("I4").Value = 355
$ws.Range("K4").Value = 1065
$ws.Range("M4").Value = -953
$ws.Range("H23").Value = 344.8889
$ws.Range("I23").Value = 200
$ws.Range("J23").Value = 417.33334
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 1252.00002
$ws.Range("M23").Value = -365
$ws.Range("N23").Value = -1722.00002
$ws.Range("H43").Value = 12893
$ws.Range("J43").Value = 13219.2
$ws.Range("L43").Value = 39657.60000000001
$ws.Range("N43").Value = -39885.60000000001
$ws.Range("H50").Value = 393.33334
$ws.Range("I50").Value = 385.75
$ws.Range("K50").Value = 1157.25
$ws.Range("M50").Value = -676.25
$ws.Range("H53").Value = 393.33334
$ws.Range("I53").Value = 385.75
$ws.Range("K53").Value = 1157.25
$ws.Range("M53").Value = -676.25
$ws.Range("H64").Value = 12793.857
$ws.Range("I64").Value = 7131.4
$ws.Range("J64").Value = 26950
$ws.Range("K64").Value = 21394.2
$ws.Range("L64").Value = 80850
$ws.Range("M64").Value = -21124.2
$ws.Range("N64").Value = -81390
$ws.Range("H67").Value = 12793.857
$ws.Range("I67").Value = 7131.4
$ws.Range("J67").Value = 26950
$ws.Range("K67").Value = 21394.2
$ws.Range("L67").Value = 80850
$ws.Range("M67").Value = -20458.2
$ws.Range("N67").Value = -82722
$ws.Range("H103").Value = 418.625
$ws.Range("I103").Value = 335.57144
$ws.Range("K103").Value = 1006.71432
$ws.Range("M103").Value = -127.71432
$ws.Range("H113").Value = 390.5357
$ws.Range("J113").Value = 500.3889
$ws.Range("L113").Value = 1501.1667
$ws.Range("N113").Value = -5841.1667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 1995
$ws.Range("J21").Value = 1995
$ws.Range("L21").Value = 1995
$ws.Range("N21").Value = -2341
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H30").Value = 1995
$ws.Range("J30").Value = 1995
$ws.Range("L30").Value = 1995
$ws.Range("N30").Value = -2205
$ws.Range("H102").Value = 915.8
$ws.Range("I102").Value = 905.0909
$ws.Range("J102").Value = 994.3333
$ws.Range("K102").Value = 905.0909
$ws.Range("L102").Value = 994.3333
$ws.Range("M102").Value = 716.9091
$ws.Range("N102").Value = -4238.3333
$ws.Range("H107").Value = 439.66666
$ws.Range("I107").Value = 280.16666
$ws.Range("J107").Value = 758.6667
$ws.Range("K107").Value = 280.16666
$ws.Range("L107").Value = 758.6667
$ws.Range("M107").Value = 1639.83334
$ws.Range("N107").Value = -4598.6667
$ws.Range("H113").Value = 4125
$ws.Range("I113").Value = 4743.5
$ws.Range("J113").Value = 3506.5
$ws.Range("K113").Value = 4743.5
$ws.Range("L113").Value = 3506.5
$ws.Range("M113").Value = -2573.5
$ws.Range("N113").Value = -7846.5
$ws.Range("H136").Value = 39197.26
$ws.Range("J136").Value = 39197.26
$ws.Range("L136").Value = 117591.78
$ws.Range("N136").Value = -122691.78

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1167.8572
$ws.Range("I35").Value = 1167.8572
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1167.8572
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -831.8571999999999
$ws.Range("H44").Value = 6939
$ws.Range("J44").Value = 6939
$ws.Range("L44").Value = 6939
$ws.Range("N44").Value = -7851
$ws.Range("H46").Value = 4681.357
$ws.Range("I46").Value = 2001
$ws.Range("J46").Value = 4887.5386
$ws.Range("K46").Value = 2001
$ws.Range("L46").Value = 4887.5386
$ws.Range("M46").Value = -1813
$ws.Range("N46").Value = -5263.5386
$ws.Range("H61").Value = 1913.6471
$ws.Range("I61").Value = 1532.9375
$ws.Range("J61").Value = 8005
$ws.Range("K61").Value = 1532.9375
$ws.Range("L61").Value = 8005
$ws.Range("M61").Value = -1330.9375
$ws.Range("N61").Value = -8409
$ws.Range("H68").Value = 3985.2856
$ws.Range("J68").Value = 4949.5
$ws.Range("L68").Value = 4949.5
$ws.Range("N68").Value = -6447.5
$ws.Range("H71").Value = 3985.2856
$ws.Range("J71").Value = 4949.5
$ws.Range("L71").Value = 24747.5
$ws.Range("N71").Value = -32235.5
$ws.Range("H86").Value = 70000
$ws.Range("J86").Value = 70000
$ws.Range("L86").Value = 70000
$ws.Range("N86").Value = -72372
$ws.Range("H89").Value = 70000
$ws.Range("J89").Value = 70000
$ws.Range("L89").Value = 210000
$ws.Range("N89").Value = -221856
$ws.Range("H100").Value = 380499.5
$ws.Range("I100").Value = 433427.44
$ws.Range("K100").Value = 433427.44
$ws.Range("M100").Value = -432886.44
$ws.Range("H113").Value = 1913.6471
$ws.Range("I113").Value = 1532.9375
$ws.Range("J113").Value = 8005
$ws.Range("K113").Value = 1532.9375
$ws.Range("L113").Value = 8005
$ws.Range("M113").Value = 637.0625
$ws.Range("N113").Value = -12345
$ws.Range("H122").Value = 5187.25
$ws.Range("I122").Value = 3546.4
$ws.Range("K122").Value = 10639.2
$ws.Range("M122").Value = -8189.200000000001
$ws.Range("H136").Value = 28575552
$ws.Range("I136").Value = 3636.3704
$ws.Range("K136").Value = 10909.1112
$ws.Range("M136").Value = -8359.111199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 39998.61
$ws.Range("I15").Value = 39998
$ws.Range("K15").Value = 39998
$ws.Range("M15").Value = -39710
$ws.Range("H32").Value = 7911.25
$ws.Range("I32").Value = 7911.25
$ws.Range("K32").Value = 7911.25
$ws.Range("M32").Value = -7594.25
$ws.Range("H34").Value = 27559
$ws.Range("I34").Value = 26618.334
$ws.Range("K34").Value = 26618.334
$ws.Range("M34").Value = -26415.334
$ws.Range("H62").Value = 4884.615
$ws.Range("J62").Value = 4625
$ws.Range("L62").Value = 4625
$ws.Range("N62").Value = -5873
$ws.Range("H65").Value = 4884.615
$ws.Range("J65").Value = 4625
$ws.Range("L65").Value = 23125
$ws.Range("N65").Value = -29365
$ws.Range("H88").Value = 37118.668
$ws.Range("J88").Value = 37118.668
$ws.Range("L88").Value = 37118.668
$ws.Range("N88").Value = -37930.668
$ws.Range("H91").Value = 37118.668
$ws.Range("J91").Value = 37118.668
$ws.Range("L91").Value = 37118.668
$ws.Range("N91").Value = -39926.668
$ws.Range("H122").Value = 3789.7778
$ws.Range("I122").Value = 3789.7778
$ws.Range("K122").Value = 11369.3334
$ws.Range("M122").Value = -8919.3334
$ws.Range("H126").Value = 2655.889
$ws.Range("I126").Value = 2793.9167
$ws.Range("K126").Value = 8381.750100000001
$ws.Range("M126").Value = -5911.750100000001
$ws.Range("H132").Value = 4444.8184
$ws.Range("I132").Value = 4444.8184
$ws.Range("K132").Value = 13334.4552
$ws.Range("M132").Value = -10804.4552
$ws.Range("H136").Value = 1836.8096
$ws.Range("I136").Value = 1836.8096
$ws.Range("K136").Value = 5510.4288
$ws.Range("M136").Value = -2960.4288
